$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '58.083.73'
$ws.Cells.Item(2, 5).Value = '  -1.05%  '

$ws.Cells.Item(3, 4).Value = '3.117.31'
$ws.Cells.Item(3, 5).Value = '  +0.98%  '

$ws.Cells.Item(4, 5).Value = '  +0.00%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '527.59'
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).Value = '  +1.16%  '

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '142.26'
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).Value = '  -1.23%  '

$ws.Cells.Item(7, 5).Value = '  +0.08%  '

$ws.Cells.Item(8, 4).Value = '3.118.18'
$ws.Cells.Item(8, 5).Value = '  +1.01%  '

$ws.Cells.Item(9, 5).Value = '  +1.68%  '

$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '7.16'
$ws.Cells.Item(10, 4).ClearFormats()
$ws.Cells.Item(10, 5).Value = '  -2.56%  '

$ws.Cells.Item(11, 5).Value = '  -0.46%  '

$ws.Cells.Item(12, 5).Value = '  +2.36%  '

$ws.Cells.Item(13, 4).Value = '3.656.79'

$ws.Cells.Item(14, 5).Value = '  +3.19%  '

$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '25.63'
$ws.Cells.Item(15, 4).ClearFormats()
$ws.Cells.Item(15, 5).Value = '  -4.11%  '

$ws.Cells.Item(16, 5).Value = '  -1.06%  '

$ws.Cells.Item(17, 4).Value = '58.141.11'
$ws.Cells.Item(17, 5).Value = '  -0.96%  '

$ws.Cells.Item(18, 4).Value = '3.107.61'
$ws.Cells.Item(18, 5).Value = '  +0.68%  '

$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '6.13'
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Cells.Item(19, 5).Value = '  -0.52%  '

$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '12.77'
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(20, 5).Value = '  -0.98%  '

$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '7.97'
$ws.Cells.Item(21, 4).ClearFormats()
$ws.Cells.Item(21, 5).Value = '  -1.91%  '

$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '342.35'
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(22, 5).Value = '  -0.04%  '

$ws.Cells.Item(24, 5).Value = '  +1.73%  '

$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '67.61'
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(25, 5).Value = '  +2.96%  '

$ws.Cells.Item(26, 5).Value = '  -0.96%  '

$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '1.00'
$ws.Cells.Item(27, 4).ClearFormats()
$ws.Cells.Item(27, 5).Value = '  +0.17%  '

$ws.Cells.Item(28, 4).Value = '0.0₃0922'
$ws.Cells.Item(28, 5).Value = '  +0.01%  '

$ws.Cells.Item(29, 5).Value = '  +0.02%  '

$ws.Cells.Item(30, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '7.32'
$ws.Cells.Item(30, 4).ClearFormats()
$ws.Cells.Item(30, 5).Value = '  +0.84%  '

$ws.Cells.Item(31, 2).Value = 'RenderToken'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '6.38'
$ws.Cells.Item(31, 4).ClearFormats()
$ws.Cells.Item(31, 5).Value = '  -3.69%  '

$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '1.87'
$ws.Cells.Item(32, 4).ClearFormats()
$ws.Cells.Item(32, 5).Value = '  +1.68%  '

$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '21.04'
$ws.Cells.Item(33, 4).ClearFormats()
$ws.Cells.Item(33, 5).Value = '  +0.17%  '

$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '1.18'
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(34, 5).Value = '  -1.08%  '

$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '158.39'
$ws.Cells.Item(35, 4).ClearFormats()
$ws.Cells.Item(35, 5).Value = '  +2.77%  '

$ws.Cells.Item(36, 5).Value = '  +2.16%  '

$ws.Cells.Item(37, 5).Value = '  +1.25%  '

$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '26.38'
$ws.Cells.Item(38, 4).ClearFormats()
$ws.Cells.Item(38, 5).Value = '  -2.16%  '

$ws.Cells.Item(40, 5).Value = '  +13.39%  '

$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.0666'
$ws.Cells.Item(41, 4).ClearFormats()
$ws.Cells.Item(41, 5).Value = '  -2.88%  '

$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '3.98'
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(42, 5).Value = '  +1.82%  '

$ws.Cells.Item(43, 5).Value = '  +4.39%  '

$ws.Cells.Item(44, 4).Value = '3.160.45'
$ws.Cells.Item(44, 5).Value = '  +1.03%  '

$ws.Cells.Item(45, 5).Value = '  -0.39%  '

$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '1.00'
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(46, 5).Value = '  -0.03%  '

$ws.Cells.Item(47, 5).Value = '  +2.89%  '

$ws.Cells.Item(48, 4).Value = '2.276.71'
$ws.Cells.Item(48, 5).Value = '  -0.32%  '

$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '0.998'
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Cells.Item(49, 5).Value = '  +4.14%  '

$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '20.63'
$ws.Cells.Item(51, 4).ClearFormats()
$ws.Cells.Item(51, 5).Value = '  -0.07%  '
